# "Generate Report for Handback"
#
# For both the zh-cn and de-de localization-status worksheets, the
# handoff rows (row 2 = 5b1d5279-..., row 3 = 83fb00ad-...) are updated
# to reflect that the handback has happened:
#   - Status (col B) moves from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - Latest Target File (col E) / Latest Handback File (col F) are
#     populated (they mirror the Source File Name / Latest Handoff File
#     that was handed back)
#   - Latest Handback DateTime (col G) is stamped with the handback time

$wb = $excel.ActiveWorkbook

$mdUrlPrefix  = "https://github.com/OpenLocalizationTest/oltest/blob/fbc13ee1a41ee0b45c7e4916ffe4a4d1c3f0165b/e2e/"
$mdFile1      = "5b1d5279-bf6b-45b9-ab9b-eb49fc76bace.md"
$mdFile2      = "83fb00ad-0850-4934-9fa6-8d99ff61d73b.md"

$handedBackStatus = "Handed back: in sync with en-US"

function Update-LocalizationSheet {
    param($SheetName, $XlfFile1, $XlfFile2, $XlfUrlPrefix, $HandbackTime2, $HandbackTime3)

    $ws = $wb.Worksheets.Item($SheetName)

    # --- Row 2 (5b1d5279-...) ---
    $ws.Range("B2").Value = $handedBackStatus

    $ws.Range("E2").Value = $mdFile1
    $ws.Hyperlinks.Add($ws.Range("E2"), ($mdUrlPrefix + $mdFile1), [Type]::Missing, [Type]::Missing, $mdFile1)

    $ws.Range("F2").Value = $XlfFile1
    $ws.Hyperlinks.Add($ws.Range("F2"), ($XlfUrlPrefix + $XlfFile1), [Type]::Missing, [Type]::Missing, $XlfFile1)

    $ws.Range("G2").Value = $HandbackTime2

    # --- Row 3 (83fb00ad-...) ---
    $ws.Range("B3").Value = $handedBackStatus

    $ws.Range("E3").Value = $mdFile2
    $ws.Hyperlinks.Add($ws.Range("E3"), ($mdUrlPrefix + $mdFile2), [Type]::Missing, [Type]::Missing, $mdFile2)

    $ws.Range("F3").Value = $XlfFile2
    $ws.Hyperlinks.Add($ws.Range("F3"), ($XlfUrlPrefix + $XlfFile2), [Type]::Missing, [Type]::Missing, $XlfFile2)

    $ws.Range("G3").Value = $HandbackTime3
}

# zh-cn sheet
Update-LocalizationSheet "zh-cn" `
    "5b1d5279-bf6b-45b9-ab9b-eb49fc76bace.6df0435f002135aba3585017242ba133463ce284.zh-cn.xlf" `
    "83fb00ad-0850-4934-9fa6-8d99ff61d73b.f7e97f9cb934a04dade58fdb7070723b9d28e2d5.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86c4005fbe442adbcc0d4dd47b689349991806bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/" `
    "2016-03-01 09:46:58" `
    "2016-03-01 09:46:58"

# de-de sheet
Update-LocalizationSheet "de-de" `
    "5b1d5279-bf6b-45b9-ab9b-eb49fc76bace.6df0435f002135aba3585017242ba133463ce284.de-de.xlf" `
    "83fb00ad-0850-4934-9fa6-8d99ff61d73b.f7e97f9cb934a04dade58fdb7070723b9d28e2d5.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2a07a5fcd472e98ae093f9378c5cb911485ac271/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/" `
    "2016-03-01 09:47:17" `
    "2016-03-01 09:47:17"
